# Update EC (Estado de Cuenta) database:
# 1) The "Periodo Mora" (arrears period) labels in column E, rows 16-29,
#    are re-ordered from descending (2105 -> 2004) to ascending (2004 -> 2105).
# 2) The "Valor Mora" amounts in F16 and F29 are swapped to stay aligned
#    with their (now re-ordered) period rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the Valor Mora values for the first and last period rows ---
$ws.Range("F16").Value = 26919
$ws.Range("F29").Value = 28090

# --- Re-order the Periodo Mora labels (column E) to ascending order ---
$ws.Range("E16").Value = "2004"
$ws.Range("E17").Value = "2005"
$ws.Range("E18").Value = "2006"
$ws.Range("E19").Value = "2007"
$ws.Range("E20").Value = "2008"
$ws.Range("E21").Value = "2009"
$ws.Range("E22").Value = "2010"
$ws.Range("E23").Value = "2011"
$ws.Range("E24").Value = "2012"
$ws.Range("E25").Value = "2101"
$ws.Range("E26").Value = "2102"
$ws.Range("E27").Value = "2103"
$ws.Range("E28").Value = "2104"
$ws.Range("E29").Value = "2105"
